# Actualizacion automatica del mapa (2025-10-29 11:16:23)
#
# Two records (rows 25 and 26 of sheet "INCO") were resolved/removed
# upstream. Deleting the entire rows shifts every subsequent record up by
# two positions, which is exactly what the published diff shows (old row
# 27 becomes new row 25, ... old row 37 disappears because nothing is left
# to shift into it), and shrinks the used range from A1:R37 to A1:R35.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("25:26").Delete()
